$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the three formulas so the inner SUBSTITUTE replaces "\" with "\" (single backslash)
# instead of "\\" (double backslash), which removes the doubled backslashes from the results.
$ws.Range("C2").Formula = '=SUBSTITUTE(SUBSTITUTE(B2,"W:","\\val-fs01\EOLE-SAP"),"\","\")'
$ws.Range("C3").Formula = '=SUBSTITUTE(SUBSTITUTE(B3,"S:","\\val-fs01\Services"),"\","\")'
$ws.Range("C4").Formula = '=SUBSTITUTE(SUBSTITUTE(B4,"S:","\\val-fs01\Services"),"\","\")'

# Widen column C slightly (target stored width is 62.28515625 characters;
# the engine quantizes ColumnWidth to pixel-based steps the same way Excel does,
# so 61.5 is the input that rounds to the closest achievable stored width)
$ws.Columns("C").ColumnWidth = 61.5

# Move the active selection to C2
$ws.Range("C2").Select()
